# Update "想去人数" (want-to-go count) values across sheets, matching the
# site's re-generated output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1039
$wsExpo.Range("F3").Value = 41
$wsExpo.Range("F4").Value = 2268
$wsExpo.Range("F5").Value = 19
$wsExpo.Range("F6").Value = 498

# --- Sheet "演出" (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 18

# --- Sheet "全部类型" (All types, combines the above) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 18
$wsAll.Range("F4").Value = 1039
$wsAll.Range("F5").Value = 41
$wsAll.Range("F6").Value = 2268
$wsAll.Range("F7").Value = 19
$wsAll.Range("F8").Value = 498
